$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores prices/percentages (columns D and E) as text, e.g. "304.78" or "1.14%".
# Pre-format the affected range as Text so the COM layer does not silently reinterpret
# these numeric-looking strings as numbers (which would lose trailing zeros / "%" signs).
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @(
    @{ Cell = "D2"; Value = "304.84" }
    @{ Cell = "E2"; Value = "1.09%" }
    @{ Cell = "D3"; Value = "35.88" }
    @{ Cell = "E3"; Value = "0.98%" }
    @{ Cell = "D4"; Value = "5.033" }
    @{ Cell = "E4"; Value = "-0.80%" }
    @{ Cell = "D5"; Value = "0.08084" }
    @{ Cell = "E5"; Value = "1.12%" }
    @{ Cell = "D6"; Value = "1.911" }
    @{ Cell = "E6"; Value = "-0.99%" }
    @{ Cell = "D7"; Value = "4.140" }
    @{ Cell = "E7"; Value = "2.41%" }
    @{ Cell = "D8"; Value = "7.839" }
    @{ Cell = "E8"; Value = "0.93%" }
    @{ Cell = "D9"; Value = "0.9324" }
    @{ Cell = "E9"; Value = "0.63%" }
    @{ Cell = "D10"; Value = "0.1265" }
    @{ Cell = "E10"; Value = "-17.37%" }
    @{ Cell = "D11"; Value = "0.1914" }
    @{ Cell = "E11"; Value = "0.48%" }
    @{ Cell = "D12"; Value = "0.09184" }
    @{ Cell = "E12"; Value = "2.21%" }
    @{ Cell = "D13"; Value = "0.03495" }
    @{ Cell = "E13"; Value = "1.20%" }
    @{ Cell = "D14"; Value = "0.09924" }
    @{ Cell = "E14"; Value = "0.39%" }
    @{ Cell = "D15"; Value = "0.001413" }
    @{ Cell = "E15"; Value = "1.13%" }
    @{ Cell = "D16"; Value = "0.006637" }
    @{ Cell = "E16"; Value = "15.32%" }
    @{ Cell = "D17"; Value = "3.619" }
    @{ Cell = "E17"; Value = "2.48%" }
    @{ Cell = "D18"; Value = "3.199" }
    @{ Cell = "E18"; Value = "8.28%" }
    @{ Cell = "D19"; Value = "0.3448" }
    @{ Cell = "E19"; Value = "0.13%" }
    @{ Cell = "B20"; Value = "MCDex" }
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb" }
    @{ Cell = "D20"; Value = "5.207" }
    @{ Cell = "B21"; Value = "ProBitToken" }
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob" }
    @{ Cell = "D21"; Value = "0.1306" }
    @{ Cell = "E21"; Value = "0.29%" }
    @{ Cell = "D22"; Value = "0.2533" }
    @{ Cell = "E22"; Value = "5.69%" }
    @{ Cell = "D23"; Value = "0.04410" }
    @{ Cell = "E23"; Value = "-1.86%" }
    @{ Cell = "D24"; Value = "0.001236" }
    @{ Cell = "E24"; Value = "1.98%" }
    @{ Cell = "D25"; Value = "0.004725" }
    @{ Cell = "E25"; Value = "-1.01%" }
    @{ Cell = "D26"; Value = "0.0001301" }
    @{ Cell = "E26"; Value = "5.86%" }
    @{ Cell = "D27"; Value = "0.0003132" }
    @{ Cell = "E27"; Value = "3.66%" }
    @{ Cell = "D39"; Value = "0.01958" }
    @{ Cell = "E39"; Value = "6.11%" }
    @{ Cell = "E40"; Value = "8.40%" }
    @{ Cell = "D41"; Value = "0.007595" }
    @{ Cell = "E41"; Value = "4.06%" }
    @{ Cell = "D42"; Value = "0.01017" }
    @{ Cell = "E42"; Value = "-4.03%" }
    @{ Cell = "D43"; Value = "0.1368" }
    @{ Cell = "E43"; Value = "2.87%" }
    @{ Cell = "D44"; Value = "0.002102" }
    @{ Cell = "E44"; Value = "-0.31%" }
    @{ Cell = "D45"; Value = "0.01069" }
    @{ Cell = "E45"; Value = "10.03%" }
    @{ Cell = "D46"; Value = "0.00006357" }
    @{ Cell = "E46"; Value = "2.16%" }
    @{ Cell = "D47"; Value = "0.00000000751" }
    @{ Cell = "E47"; Value = "0.20%" }
    @{ Cell = "D48"; Value = "64.96" }
    @{ Cell = "E48"; Value = "0.45%" }
    @{ Cell = "E49"; Value = "-3.48%" }
    @{ Cell = "D50"; Value = "0.00002102" }
    @{ Cell = "E50"; Value = "0.20%" }
    @{ Cell = "D51"; Value = "0.0002002" }
    @{ Cell = "E51"; Value = "0.20%" }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value2 = $u.Value
}

Write-Host "Applied $($updates.Count) cell updates"
